# CF Expense Request - 8th Jan 2025
$wb = $excel.ActiveWorkbook

# --- Approver sheet: update the hard-coded password value ---
$approverWs = $wb.Worksheets.Item("Approver")
$approverWs.Range("B2").Value = "Bingo@12345"
$approverWs.Activate()
$approverWs.Range("C11").Select()

# --- Actions sheet: reorder the approver action list ---
$actionsWs = $wb.Worksheets.Item("Actions")
$actionsWs.Range("A2").Value = "Approve"
$actionsWs.Range("A3").Value = "Reject"
$actionsWs.Range("A4").Value = "Delete"
$actionsWs.Range("A5").Value = "Request More Information"
$actionsWs.Range("A6").Value = "Edit"
$actionsWs.Activate()
$actionsWs.Range("F13").Select()
